$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.809.49"
$ws.Range("E2").Value = "  -5.35%  "
$ws.Range("D3").Value = "2.652.11"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.588"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.52"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0821"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.90"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("D13").Value = "3.061.55"
$ws.Range("E13").Value = "  +2.30%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "2.656.81"
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.900"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "43.715.48"
$ws.Range("E18").Value = "  -5.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "0.0₃0989"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.79"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "273.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  +4.08%  "
$ws.Range("E25").Value = "  -1.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.57%  "
$ws.Range("E29").Value = "  -3.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.07%  "
$ws.Range("E39").Value = "  -0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0320"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "2.104.95"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("D48").Value = "2.913.21"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.56"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("E51").Value = "  -2.81%  "
